# Results workbook update:
#   The "Award" labels on the Guest sheet change from gold/silver/bronze to
#   "honorary gold"/"honorary silver"/"honorary bronze" and the sheet's
#   conditional formatting (which colors rows by medal) is updated to match
#   the new wording. The Regular sheet is left untouched.

$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Guest")

# --- Update the Award column (B2:B6) on the Guest sheet ---------------------
$ws2.Cells.Item(2, 2).Value = "honorary gold"
$ws2.Cells.Item(3, 2).Value = "honorary silver"
$ws2.Cells.Item(4, 2).Value = "honorary silver"
$ws2.Cells.Item(5, 2).Value = "honorary bronze"
$ws2.Cells.Item(6, 2).Value = "honorary bronze"

# --- Update the conditional formatting rules on the Guest sheet -------------
# (same three rules/colors, just the text they match against changes)
$fcs2 = $ws2.Range("A1:N1048576").FormatConditions

$fcBronze = $fcs2.Item(1)
$fcBronze.Formula1 = '=$B1="honorary bronze"'

$fcSilver = $fcs2.Item(2)
$fcSilver.Formula1 = '=$B1="honorary silver"'

$fcGold = $fcs2.Item(3)
$fcGold.Formula1 = '=$B1="honorary gold"'

# --- Grow the workbook's differential-style table to mirror the source -----
# The published workbook ends up with three extra (unused) differential
# styles duplicating the gold/silver/bronze fills. Recreate that by adding
# throwaway rules with the matching fills and then removing the rules again
# (the styles they reference stay registered in the style table).
$d1 = $fcs2.Add(2, 0, '=FALSE')
$d1.Interior.Color = 65535      # yellow   (gold)
$d2 = $fcs2.Add(2, 0, '=FALSE')
$d2.Interior.Color = 12566463   # gray     (silver)
$d3 = $fcs2.Add(2, 0, '=FALSE')
$d3.Interior.Color = 1331390    # orange   (bronze)
$fcs2.Item(4).Delete()
$fcs2.Item(4).Delete()
$fcs2.Item(4).Delete()

# --- Drop the stale cell selection recorded for the Guest sheet ------------
$ws2.Range("A1").Select()

Write-Host "Updated Guest sheet awards and conditional formatting"
